$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cell values: Price (D), Volume 1h (E), and Hora (G) columns
# for the refreshed crypto symbol snapshot (row => column letter => new text value).
$updates = @(
    @{ Cell = "D2"; Value = "335.79" }
    @{ Cell = "E2"; Value = "1.74%" }
    @{ Cell = "G2"; Value = "11" }
    @{ Cell = "D3"; Value = "43.91" }
    @{ Cell = "G3"; Value = "11" }
    @{ Cell = "D4"; Value = "5.757" }
    @{ Cell = "E4"; Value = "2.04%" }
    @{ Cell = "G4"; Value = "11" }
    @{ Cell = "D5"; Value = "0.08385" }
    @{ Cell = "E5"; Value = "2.15%" }
    @{ Cell = "G5"; Value = "11" }
    @{ Cell = "D6"; Value = "8.837" }
    @{ Cell = "E6"; Value = "0.87%" }
    @{ Cell = "G6"; Value = "11" }
    @{ Cell = "D7"; Value = "1.959" }
    @{ Cell = "E7"; Value = "-2.63%" }
    @{ Cell = "G7"; Value = "11" }
    @{ Cell = "D8"; Value = "2.879" }
    @{ Cell = "E8"; Value = "-2.89%" }
    @{ Cell = "G8"; Value = "11" }
    @{ Cell = "D9"; Value = "0.9520" }
    @{ Cell = "E9"; Value = "3.34%" }
    @{ Cell = "G9"; Value = "11" }
    @{ Cell = "D10"; Value = "0.1249" }
    @{ Cell = "E10"; Value = "-2.11%" }
    @{ Cell = "G10"; Value = "11" }
    @{ Cell = "D11"; Value = "0.1982" }
    @{ Cell = "E11"; Value = "1.26%" }
    @{ Cell = "G11"; Value = "11" }
    @{ Cell = "D12"; Value = "0.1017" }
    @{ Cell = "E12"; Value = "8.11%" }
    @{ Cell = "G12"; Value = "11" }
    @{ Cell = "D13"; Value = "0.04465" }
    @{ Cell = "E13"; Value = "16.31%" }
    @{ Cell = "G13"; Value = "11" }
    @{ Cell = "D14"; Value = "0.1068" }
    @{ Cell = "E14"; Value = "0.69%" }
    @{ Cell = "G14"; Value = "11" }
    @{ Cell = "D15"; Value = "0.001292" }
    @{ Cell = "E15"; Value = "-1.17%" }
    @{ Cell = "G15"; Value = "11" }
    @{ Cell = "D16"; Value = "0.005938" }
    @{ Cell = "E16"; Value = "-3.04%" }
    @{ Cell = "G16"; Value = "11" }
    @{ Cell = "D17"; Value = "3.491" }
    @{ Cell = "E17"; Value = "1.29%" }
    @{ Cell = "G17"; Value = "11" }
    @{ Cell = "D18"; Value = "4.513" }
    @{ Cell = "E18"; Value = "-0.07%" }
    @{ Cell = "G18"; Value = "11" }
    @{ Cell = "G19"; Value = "11" }
    @{ Cell = "D20"; Value = "8.715" }
    @{ Cell = "E20"; Value = "4.99%" }
    @{ Cell = "G20"; Value = "11" }
    @{ Cell = "D21"; Value = "0.1363" }
    @{ Cell = "E21"; Value = "-0.79%" }
    @{ Cell = "G21"; Value = "11" }
    @{ Cell = "E22"; Value = "-0.79%" }
    @{ Cell = "G22"; Value = "11" }
    @{ Cell = "D23"; Value = "0.04414" }
    @{ Cell = "E23"; Value = "0.17%" }
    @{ Cell = "G23"; Value = "11" }
    @{ Cell = "E24"; Value = "0.15%" }
    @{ Cell = "G24"; Value = "11" }
    @{ Cell = "D25"; Value = "0.004362" }
    @{ Cell = "E25"; Value = "1.09%" }
    @{ Cell = "G25"; Value = "11" }
    @{ Cell = "D26"; Value = "0.0001263" }
    @{ Cell = "E26"; Value = "5.19%" }
    @{ Cell = "G26"; Value = "11" }
    @{ Cell = "D27"; Value = "0.0003994" }
    @{ Cell = "E27"; Value = "-94.68%" }
    @{ Cell = "G27"; Value = "11" }
    @{ Cell = "G28"; Value = "11" }
    @{ Cell = "G29"; Value = "11" }
    @{ Cell = "G30"; Value = "11" }
    @{ Cell = "G31"; Value = "11" }
    @{ Cell = "G32"; Value = "11" }
    @{ Cell = "G33"; Value = "11" }
    @{ Cell = "G34"; Value = "11" }
    @{ Cell = "G35"; Value = "11" }
    @{ Cell = "G36"; Value = "11" }
    @{ Cell = "G37"; Value = "11" }
    @{ Cell = "G38"; Value = "11" }
    @{ Cell = "D39"; Value = "0.02814" }
    @{ Cell = "E39"; Value = "1.92%" }
    @{ Cell = "G39"; Value = "11" }
    @{ Cell = "D40"; Value = "0.06036" }
    @{ Cell = "E40"; Value = "9.50%" }
    @{ Cell = "G40"; Value = "11" }
    @{ Cell = "D41"; Value = "0.007931" }
    @{ Cell = "E41"; Value = "0.02%" }
    @{ Cell = "G41"; Value = "11" }
    @{ Cell = "D42"; Value = "0.1428" }
    @{ Cell = "E42"; Value = "0.61%" }
    @{ Cell = "G42"; Value = "11" }
    @{ Cell = "D43"; Value = "0.008970" }
    @{ Cell = "E43"; Value = "0.21%" }
    @{ Cell = "G43"; Value = "11" }
    @{ Cell = "D44"; Value = "0.002147" }
    @{ Cell = "E44"; Value = "-1.12%" }
    @{ Cell = "G44"; Value = "11" }
    @{ Cell = "D45"; Value = "0.01016" }
    @{ Cell = "E45"; Value = "-10.88%" }
    @{ Cell = "G45"; Value = "11" }
    @{ Cell = "D46"; Value = "0.00007328" }
    @{ Cell = "E46"; Value = "8.08%" }
    @{ Cell = "G46"; Value = "11" }
    @{ Cell = "E47"; Value = "0.18%" }
    @{ Cell = "G47"; Value = "11" }
    @{ Cell = "D48"; Value = "0.003202" }
    @{ Cell = "E48"; Value = "0.36%" }
    @{ Cell = "G48"; Value = "11" }
    @{ Cell = "D49"; Value = "0.002273" }
    @{ Cell = "E49"; Value = "-0.32%" }
    @{ Cell = "G49"; Value = "11" }
    @{ Cell = "E50"; Value = "0.18%" }
    @{ Cell = "G50"; Value = "11" }
    @{ Cell = "E51"; Value = "0.18%" }
    @{ Cell = "G51"; Value = "11" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so values like "335.79" / "1.74%" / "11" stay as literal
    # text (matching the source data format) instead of being auto-converted to
    # numbers or percentages by Excel's smart-entry parsing.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
